# "Add files via upload" — adds a new "Emisiones por Sectores" worksheet
# (sector-level emissions breakdown) after the existing "Emisiones" sheet,
# and normalizes a couple of cosmetic style/number-format bookkeeping
# details left over from the previous edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new worksheet right after "Emisiones" (the last sheet),
#    so it lands at the end of the tab strip and becomes the active tab,
#    exactly like Excel does when you ctrl-click "New Sheet" at the end.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Emisiones por Sectores"

# ---------------------------------------------------------------------
# 2. Header row 1: hex color codes for each sector, columns B..I
# ---------------------------------------------------------------------
$ws3.Cells.Item(1, 2).Value = "#244b7e"
$ws3.Cells.Item(1, 3).Value = "#48a04a"
$ws3.Cells.Item(1, 4).Value = "#5f5ca5"
$ws3.Cells.Item(1, 5).Value = "#8d0049"
$ws3.Cells.Item(1, 6).Value = "#f7bd0b"
$ws3.Cells.Item(1, 7).Value = "#066e67"
$ws3.Cells.Item(1, 8).Value = "#998b00"
$ws3.Cells.Item(1, 9).Value = "#8a3702"

# ---------------------------------------------------------------------
# 3. Header row 2: "Año" plus sector names, columns A..I
# ---------------------------------------------------------------------
$ws3.Cells.Item(2, 1).Value = "Año"
$ws3.Cells.Item(2, 2).Value = "Industrias manufactureras y de la construcción"
$ws3.Cells.Item(2, 3).Value = "Agropecuario"
$ws3.Cells.Item(2, 4).Value = "Comercial"
$ws3.Cells.Item(2, 5).Value = "Forestal"
$ws3.Cells.Item(2, 6).Value = "Industrias de la energía"
$ws3.Cells.Item(2, 7).Value = "Residencial"
$ws3.Cells.Item(2, 8).Value = "Saneamiento"
$ws3.Cells.Item(2, 9).Value = "Transporte"

# ---------------------------------------------------------------------
# 4. Data rows 3..7: Año + per-sector Mt CO2eq values, years 2010-2014
# ---------------------------------------------------------------------
$ws3.Cells.Item(3, 1).Value = 2010
$ws3.Cells.Item(3, 2).Value = 25.103353904482901
$ws3.Cells.Item(3, 3).Value = 54.696658991737003
$ws3.Cells.Item(3, 4).Value = 1.2313263955931999
$ws3.Cells.Item(3, 5).Value = 84.215286458579996
$ws3.Cells.Item(3, 6).Value = 26.664509108720001
$ws3.Cells.Item(3, 7).Value = 4.7085485915817999
$ws3.Cells.Item(3, 8).Value = 9.2159828227689999
$ws3.Cells.Item(3, 9).Value = 24.525282009076101

$ws3.Cells.Item(4, 1).Value = 2011
$ws3.Cells.Item(4, 2).Value = 24.5130436971344
$ws3.Cells.Item(4, 3).Value = 52.521634329784
$ws3.Cells.Item(4, 4).Value = 1.3504843117316001
$ws3.Cells.Item(4, 5).Value = 88.608263763509001
$ws3.Cells.Item(4, 6).Value = 23.376848111558999
$ws3.Cells.Item(4, 7).Value = 4.9917729544035998
$ws3.Cells.Item(4, 8).Value = 8.9828380548369999
$ws3.Cells.Item(4, 9).Value = 26.230202244301399

$ws3.Cells.Item(5, 1).Value = 2012
$ws3.Cells.Item(5, 2).Value = 26.265865562873898
$ws3.Cells.Item(5, 3).Value = 53.239591457416999
$ws3.Cells.Item(5, 4).Value = 1.6437233154972
$ws3.Cells.Item(5, 5).Value = 87.646347339390005
$ws3.Cells.Item(5, 6).Value = 24.872885487358001
$ws3.Cells.Item(5, 7).Value = 4.79309715464065
$ws3.Cells.Item(5, 8).Value = 9.2237567816109998
$ws3.Cells.Item(5, 9).Value = 27.131972344652201

$ws3.Cells.Item(6, 1).Value = 2013
$ws3.Cells.Item(6, 2).Value = 26.529880456496301
$ws3.Cells.Item(6, 3).Value = 53.492224244860999
$ws3.Cells.Item(6, 4).Value = 1.8121416057048001
$ws3.Cells.Item(6, 5).Value = 76.510767123190007
$ws3.Cells.Item(6, 6).Value = 31.814383668801
$ws3.Cells.Item(6, 7).Value = 4.7486904168702999
$ws3.Cells.Item(6, 8).Value = 9.5734717976759995
$ws3.Cells.Item(6, 9).Value = 27.295577064068599

$ws3.Cells.Item(7, 1).Value = 2014
$ws3.Cells.Item(7, 2).Value = 27.631099662162001
$ws3.Cells.Item(7, 3).Value = 52.001282318775999
$ws3.Cells.Item(7, 4).Value = 2.074727095154
$ws3.Cells.Item(7, 5).Value = 79.227719760710002
$ws3.Cells.Item(7, 6).Value = 32.322009852826
$ws3.Cells.Item(7, 7).Value = 4.9391703338299999
$ws3.Cells.Item(7, 8).Value = 9.8195287522779999
$ws3.Cells.Item(7, 9).Value = 28.957844589055

# ---------------------------------------------------------------------
# 5. Selection on the new sheet matches the author's last cursor spot.
# ---------------------------------------------------------------------
$ws3.Range("D9").Select()
